# Auto-generated edit script: updates crypto price/volume table cells
# per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.922.75"
$ws.Range("E2").Value = "  +3.20%  "
$ws.Range("D3").Value = "3.096.17"
$ws.Range("E3").Value = "  +5.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.75"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.00"
$ws.Range("E6").Value = "  +6.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.092.00"
$ws.Range("E8").Value = "  +5.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("E12").Value = "  +6.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.78"
$ws.Range("E14").Value = "  +8.72%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "3.614.68"
$ws.Range("E16").Value = "  +5.51%  "
$ws.Range("D17").Value = "66.926.18"
$ws.Range("E17").Value = "  +2.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.23"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("D19").Value = "3.102.28"
$ws.Range("E19").Value = "  +5.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.33"
$ws.Range("E20").Value = "  +18.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "469.38"
$ws.Range("E21").Value = "  +5.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("E22").Value = "  +5.84%  "
$ws.Range("E23").Value = "  +5.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.40"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +9.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.86"
$ws.Range("E26").Value = "  +7.53%  "
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.07"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("E30").Value = "  +4.35%  "
$ws.Range("E31").Value = "  +4.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000103"
$ws.Range("E32").Value = "  +4.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.26"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("E34").Value = "  +5.55%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.91"
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.84"
$ws.Range("E38").Value = "  +8.52%  "
$ws.Range("E39").Value = "  +6.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.28"
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.316"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("E42").Value = "  +4.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.73"
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.85"
$ws.Range("E44").Value = "  +3.01%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0364"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "389.77"
$ws.Range("E46").Value = "  +2.68%  "
$ws.Range("D47").Value = "2.759.78"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.86"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.80"
$ws.Range("E50").Value = "  +7.79%  "
$ws.Range("E51").Value = "  +5.46%  "
